# Update res_bus/vm_pu.xlsx results for "case with 380 kV done":
# rows 2-25 (bus indices 0-23), columns B-F and I-N get new voltage
# magnitude values from the re-run power flow case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$leftCols = @("B","C","D","E","F")
$rightCols = @("I","J","K","L","M","N")

$leftData = @(
    @(1.02,1.031784646204243,1.040412556733572,1.031289603000241,1.04911141044257),
    @(1.02,1.033058459206594,1.041440843722762,1.032380963537051,1.050377505917677),
    @(1.02,1.033882029603688,1.042105272459885,1.033086890060353,1.051196240970176),
    @(1.02,1.0342281010141,1.042384374643964,1.033383602072278,1.051540316805171),
    @(1.02,1.034286198779911,1.042431224065124,1.033433417892446,1.051598081665723),
    @(1.02,1.033886654441986,1.042109002713552,1.033090854972142,1.051200838998221),
    @(1.02,1.03221527832803,1.040760267203846,1.03165848774321,1.049539401158905),
    @(1.02,1.029264809039896,1.038376329421052,1.029132418923966,1.046607689575604),
    @(1.02,1.027294065652775,1.036782017040464,1.027446864940914,1.044650332442429),
    @(1.02,1.026439772568014,1.036090444372625,1.02671661473893,1.043802054287517),
    @(1.02,1.026122303501205,1.035833377020443,1.02644530547663,1.043486853079739),
    @(1.02,1.026190408393092,1.035888527317831,1.026503505073901,1.043554469972686),
    @(1.02,1.026413533489982,1.036069198923278,1.026694189506933,1.043776001980394),
    @(1.02,1.026550988679564,1.03618049187955,1.026811668209345,1.043912480097671),
    @(1.02,1.027350741942562,1.03682788842558,1.027495320715011,1.044706614215391),
    @(1.02,1.027852149273102,1.037233653517743,1.027924050501159,1.045204555896975),
    @(1.02,1.028144520479634,1.037470211414755,1.028174083474398,1.045494926686089),
    @(1.02,1.028244196159297,1.037550851600035,1.028259331913453,1.045593923818912),
    @(1.02,1.027798362476902,1.037190130989228,1.027878055788225,1.04515113875025),
    @(1.02,1.026347832789677,1.03601600080586,1.02663803942166,1.043710769469884),
    @(1.02,1.02543497709304,1.035276698467644,1.025858032928189,1.042804497088994),
    @(1.02,1.025918981131673,1.035668719826848,1.026271564027917,1.04328499242905),
    @(1.02,1.02782266669652,1.037209797311081,1.027898838936116,1.045175275876272),
    @(1.02,1.030028224276297,1.03899351074733,1.029785724738314,1.047366103403679),
)

$rightData = @(
    @(1.034964561482491,1.036918463248109,1.04319497772573,1.034098150403867,1.05186937860424,1.016201752684047),
    @(1.035273435708588,1.037832647183524,1.044033109662297,1.034997285372889,1.052946482001031,1.016511713292018),
    @(1.035471155117469,1.038423056106821,1.04457389496843,1.03557824410615,1.053642373762224,1.016711712602484),
    @(1.035553764025422,1.038670995666385,1.044800873609393,1.035822279470993,1.05393467360193,1.016795657426167),
    @(1.035567604400639,1.038712610130574,1.044838962820044,1.035863242425582,1.053983737255026,1.016809744238952),
    @(1.035472260954037,1.038426370139343,1.044576929310111,1.035581505701022,1.053646280478639,1.016712834806598),
    @(1.035069391078488,1.037227650998253,1.04347854865609,1.034402192658941,1.052233612815217,1.016306622863393),
    @(1.034343040084748,1.035106616685157,1.041531174743272,1.032317563598345,1.049736047611579,1.015586464119158),
    @(1.03384770550806,1.033686590967987,1.040224823884992,1.030923303807909,1.048065306455673,1.015103387114209),
    @(1.033630576519371,1.033070251106567,1.039657212668003,1.030318478007408,1.047340473282327,1.014893495769226),
    @(1.033549526622182,1.032841093166178,1.039446080928818,1.030093650784018,1.047071026044278,1.014815424423932),
    @(1.033566930150681,1.032890258404811,1.039491382819856,1.030141884627482,1.047128833028447,1.014832175908868),
    @(1.033623885045246,1.033051313381106,1.039639766497294,1.030299897159965,1.047318205029666,1.014887044582579),
    @(1.033658923996808,1.03315051518838,1.039731151331451,1.030397231563413,1.047434855200291,1.014920836600975),
    @(1.033862059812654,1.033727464444427,1.040262453036201,1.030963420678676,1.048113381686927,1.015117301761008),
    @(1.033988772461863,1.034088976913559,1.040595200159292,1.031318279388092,1.048538628907818,1.015240346950383),
    @(1.034062426590105,1.034299700392893,1.040789097678697,1.031525156296563,1.048786534382586,1.015312048117381),
    @(1.03408749748561,1.034371527879534,1.040855179860621,1.031595678039725,1.048871041095631,1.015336484681391),
    @(1.033975203799858,1.034050204612708,1.040559519030997,1.031280217407313,1.04849301779155,1.015227152510327),
    @(1.033607124257754,1.033003892850713,1.03959607938284,1.030253371064337,1.047262445568892,1.014870890123391),
    @(1.033373391797875,1.032344750115945,1.038988614888167,1.029606778861723,1.046487508386858,1.014646266510011),
    @(1.033497516743808,1.032694296833821,1.039310806334545,1.029949642674593,1.046898434586097,1.014765403471421),
    @(1.033981335682727,1.034067724580971,1.040575642378603,1.03129741630607,1.048513627905758,1.01523311472306),
    @(1.034532772410985,1.03565600350345,1.04203603744068,1.032857275917332,1.050382721911967,1.015773163242865),
)

$startRow = 2
for ($i = 0; $i -lt $leftData.Count; $i++) {
    $r = $startRow + $i
    for ($j = 0; $j -lt $leftCols.Count; $j++) {
        $col = $leftCols[$j]
        $ws.Range("$col$r").Value = $leftData[$i][$j]
    }
    for ($j = 0; $j -lt $rightCols.Count; $j++) {
        $col = $rightCols[$j]
        $ws.Range("$col$r").Value = $rightData[$i][$j]
    }
}
